$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2872560.5  # ALC!H17: 1828538.5 -> 2872560.5
$ws.Cells.Item(17, 9).Value = 0  # ALC!I17: 1500 -> 0
$ws.Cells.Item(17, 11).Value = 0  # ALC!K17: 4500 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # ALC!M17: -4332 -> (empty)
$ws.Cells.Item(29, 8).Value = 249.5  # ALC!H29: 0 -> 249.5
$ws.Cells.Item(29, 9).Value = 249  # ALC!I29: 0 -> 249
$ws.Cells.Item(29, 10).Value = 250  # ALC!J29: 0 -> 250
$ws.Cells.Item(29, 11).Value = 747  # ALC!K29: 0 -> 747
$ws.Cells.Item(29, 12).Value = 750  # ALC!L29: 0 -> 750
$ws.Cells.Item(29, 13).Value = -466  # ALC!M29: None -> -466
$ws.Cells.Item(29, 14).Value = -1312  # ALC!N29: None -> -1312
$ws.Cells.Item(40, 8).Value = 106947.17  # ALC!H40: 97030.25 -> 106947.17
$ws.Cells.Item(40, 9).Value = 501765  # ALC!I40: 430227.16 -> 501765
$ws.Cells.Item(40, 10).Value = 3951.2173  # ALC!J40: 3735.12 -> 3951.2173
$ws.Cells.Item(40, 11).Value = 501765  # ALC!K40: 430227.16 -> 501765
$ws.Cells.Item(40, 12).Value = 3951.2173  # ALC!L40: 3735.12 -> 3951.2173
$ws.Cells.Item(40, 13).Value = -501590  # ALC!M40: -430052.16 -> -501590
$ws.Cells.Item(40, 14).Value = -4301.2173  # ALC!N40: -4085.12 -> -4301.2173
$ws.Cells.Item(100, 8).Value = 4477.6665  # ALC!H100: 4517.364 -> 4477.6665
$ws.Cells.Item(100, 9).Value = 2000  # ALC!I100: 3797.3333 -> 2000
$ws.Cells.Item(100, 11).Value = 2000  # ALC!K100: 3797.3333 -> 2000
$ws.Cells.Item(100, 13).Value = -1459  # ALC!M100: -3256.3333 -> -1459
$ws.Cells.Item(111, 8).Value = 0  # ALC!H111: 1054.4 -> 0
$ws.Cells.Item(111, 9).Value = 0  # ALC!I111: 1475 -> 0
$ws.Cells.Item(111, 10).Value = 0  # ALC!J111: 774 -> 0
$ws.Cells.Item(111, 11).Value = 0  # ALC!K111: 4425 -> 0
$ws.Cells.Item(111, 12).Value = 0  # ALC!L111: 2322 -> 0
$ws.Cells.Item(111, 13).ClearContents()  # ALC!M111: -1358 -> (empty)
$ws.Cells.Item(111, 14).ClearContents()  # ALC!N111: -8456 -> (empty)
$ws.Cells.Item(132, 8).Value = 3811.6511  # ALC!H132: 4042.125 -> 3811.6511
$ws.Cells.Item(132, 9).Value = 3871.7368  # ALC!I132: 4046.8057 -> 3871.7368
$ws.Cells.Item(132, 10).Value = 3355  # ALC!J132: 4000 -> 3355
$ws.Cells.Item(132, 11).Value = 11615.2104  # ALC!K132: 12140.4171 -> 11615.2104
$ws.Cells.Item(132, 12).Value = 10065  # ALC!L132: 12000 -> 10065
$ws.Cells.Item(132, 13).Value = -9085.2104  # ALC!M132: -9610.417099999999 -> -9085.2104
$ws.Cells.Item(132, 14).Value = -15125  # ALC!N132: -17060 -> -15125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3364.353  # ARM!H45: 3512.25 -> 3364.353
$ws.Cells.Item(45, 9).Value = 2872.5454  # ARM!I45: 3060 -> 2872.5454
$ws.Cells.Item(45, 11).Value = 2872.5454  # ARM!K45: 3060 -> 2872.5454
$ws.Cells.Item(45, 13).Value = -2495.5454  # ARM!M45: -2683 -> -2495.5454
$ws.Cells.Item(61, 8).Value = 1919.3778  # ARM!H61: 1980.5416 -> 1919.3778
$ws.Cells.Item(61, 9).Value = 1929.5897  # ARM!I61: 1931.325 -> 1929.5897
$ws.Cells.Item(61, 10).Value = 1853  # ARM!J61: 2226.625 -> 1853
$ws.Cells.Item(61, 11).Value = 1929.5897  # ARM!K61: 1931.325 -> 1929.5897
$ws.Cells.Item(61, 12).Value = 1853  # ARM!L61: 2226.625 -> 1853
$ws.Cells.Item(61, 13).Value = -1717.5897  # ARM!M61: -1719.325 -> -1717.5897
$ws.Cells.Item(61, 14).Value = -2277  # ARM!N61: -2650.625 -> -2277
$ws.Cells.Item(110, 8).Value = 1469  # ARM!H110: 1649 -> 1469
$ws.Cells.Item(110, 9).Value = 831.6  # ARM!I110: 870.6667 -> 831.6
$ws.Cells.Item(110, 10).Value = 3062.5  # ARM!J110: 3050 -> 3062.5
$ws.Cells.Item(110, 11).Value = 831.6  # ARM!K110: 870.6667 -> 831.6
$ws.Cells.Item(110, 12).Value = 3062.5  # ARM!L110: 3050 -> 3062.5
$ws.Cells.Item(110, 13).Value = 1213.4  # ARM!M110: 1174.3333 -> 1213.4
$ws.Cells.Item(110, 14).Value = -7152.5  # ARM!N110: -7140 -> -7152.5
$ws.Cells.Item(128, 8).Value = 93969.664  # ARM!H128: 65995 -> 93969.664
$ws.Cells.Item(128, 10).Value = 93969.664  # ARM!J128: 65995 -> 93969.664
$ws.Cells.Item(128, 12).Value = 93969.664  # ARM!L128: 65995 -> 93969.664
$ws.Cells.Item(128, 14).Value = -103929.664  # ARM!N128: -75955 -> -103929.664
$ws.Cells.Item(132, 8).Value = 3094.9111  # ARM!H132: 3144.1555 -> 3094.9111
$ws.Cells.Item(132, 9).Value = 2823.2188  # ARM!I132: 2892.4688 -> 2823.2188
$ws.Cells.Item(132, 11).Value = 8469.6564  # ARM!K132: 8677.4064 -> 8469.6564
$ws.Cells.Item(132, 13).Value = -5939.6564  # ARM!M132: -6147.4064 -> -5939.6564
$ws.Cells.Item(135, 8).Value = 65632.42999999999  # ARM!H135: 67216 -> 65632.42999999999
$ws.Cells.Item(135, 10).Value = 65632.42999999999  # ARM!J135: 67216 -> 65632.42999999999
$ws.Cells.Item(135, 12).Value = 65632.42999999999  # ARM!L135: 67216 -> 65632.42999999999
$ws.Cells.Item(135, 14).Value = -75772.42999999999  # ARM!N135: -77356 -> -75772.42999999999
$ws.Cells.Item(136, 8).Value = 1919.3778  # ARM!H136: 1980.5416 -> 1919.3778
$ws.Cells.Item(136, 9).Value = 1929.5897  # ARM!I136: 1931.325 -> 1929.5897
$ws.Cells.Item(136, 10).Value = 1853  # ARM!J136: 2226.625 -> 1853
$ws.Cells.Item(136, 11).Value = 5788.7691  # ARM!K136: 5793.975 -> 5788.7691
$ws.Cells.Item(136, 12).Value = 5559  # ARM!L136: 6679.875 -> 5559
$ws.Cells.Item(136, 13).Value = -3238.7691  # ARM!M136: -3243.975 -> -3238.7691
$ws.Cells.Item(136, 14).Value = -10659  # ARM!N136: -11779.875 -> -10659

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 0  # BSM!H54: 11900 -> 0
$ws.Cells.Item(54, 9).Value = 0  # BSM!I54: 11900 -> 0
$ws.Cells.Item(54, 11).Value = 0  # BSM!K54: 11900 -> 0
$ws.Cells.Item(54, 13).ClearContents()  # BSM!M54: -11416 -> (empty)
$ws.Cells.Item(86, 8).Value = 2209.7778  # BSM!H86: 2220.3333 -> 2209.7778
$ws.Cells.Item(86, 9).Value = 2126.3572  # BSM!I86: 2215.7693 -> 2126.3572
$ws.Cells.Item(86, 10).Value = 2299.6155  # BSM!J86: 2224.5715 -> 2299.6155
$ws.Cells.Item(86, 11).Value = 2126.3572  # BSM!K86: 2215.7693 -> 2126.3572
$ws.Cells.Item(86, 12).Value = 2299.6155  # BSM!L86: 2224.5715 -> 2299.6155
$ws.Cells.Item(86, 13).Value = -1003.3572  # BSM!M86: -1092.7693 -> -1003.3572
$ws.Cells.Item(86, 14).Value = -4545.6155  # BSM!N86: -4470.5715 -> -4545.6155
$ws.Cells.Item(89, 8).Value = 2209.7778  # BSM!H89: 2220.3333 -> 2209.7778
$ws.Cells.Item(89, 9).Value = 2126.3572  # BSM!I89: 2215.7693 -> 2126.3572
$ws.Cells.Item(89, 10).Value = 2299.6155  # BSM!J89: 2224.5715 -> 2299.6155
$ws.Cells.Item(89, 11).Value = 10631.786  # BSM!K89: 11078.8465 -> 10631.786
$ws.Cells.Item(89, 12).Value = 11498.0775  # BSM!L89: 11122.8575 -> 11498.0775
$ws.Cells.Item(89, 13).Value = -5015.786  # BSM!M89: -5462.8465 -> -5015.786
$ws.Cells.Item(89, 14).Value = -22730.0775  # BSM!N89: -22354.8575 -> -22730.0775
$ws.Cells.Item(105, 8).Value = 3137.7273  # BSM!H105: 3811.5833 -> 3137.7273
$ws.Cells.Item(105, 9).Value = 1853  # BSM!I105: 1994.5 -> 1853
$ws.Cells.Item(105, 10).Value = 4208.3335  # BSM!J105: 4175 -> 4208.3335
$ws.Cells.Item(105, 11).Value = 1853  # BSM!K105: 1994.5 -> 1853
$ws.Cells.Item(105, 12).Value = 4208.3335  # BSM!L105: 4175 -> 4208.3335
$ws.Cells.Item(105, 13).Value = -106  # BSM!M105: -247.5 -> -106
$ws.Cells.Item(105, 14).Value = -7702.3335  # BSM!N105: -7669 -> -7702.3335
$ws.Cells.Item(107, 8).Value = 2304.3333  # BSM!H107: 1873.6666 -> 2304.3333
$ws.Cells.Item(107, 9).Value = 2450  # BSM!I107: 1621.4286 -> 2450
$ws.Cells.Item(107, 10).Value = 2013  # BSM!J107: 2756.5 -> 2013
$ws.Cells.Item(107, 11).Value = 2450  # BSM!K107: 1621.4286 -> 2450
$ws.Cells.Item(107, 12).Value = 2013  # BSM!L107: 2756.5 -> 2013
$ws.Cells.Item(107, 13).Value = -530  # BSM!M107: 298.5714 -> -530
$ws.Cells.Item(107, 14).Value = -5853  # BSM!N107: -6596.5 -> -5853

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4708.5415  # CRP!H31: 4600.08 -> 4708.5415
$ws.Cells.Item(31, 9).Value = 2570.7778  # CRP!I31: 2513.4 -> 2570.7778
$ws.Cells.Item(31, 11).Value = 2570.7778  # CRP!K31: 2513.4 -> 2570.7778
$ws.Cells.Item(31, 13).Value = -2275.7778  # CRP!M31: -2218.4 -> -2275.7778
$ws.Cells.Item(34, 8).Value = 4708.5415  # CRP!H34: 4600.08 -> 4708.5415
$ws.Cells.Item(34, 9).Value = 2570.7778  # CRP!I34: 2513.4 -> 2570.7778
$ws.Cells.Item(34, 11).Value = 2570.7778  # CRP!K34: 2513.4 -> 2570.7778
$ws.Cells.Item(34, 13).Value = -2368.7778  # CRP!M34: -2311.4 -> -2368.7778
$ws.Cells.Item(70, 8).Value = 0  # CRP!H70: 34000 -> 0
$ws.Cells.Item(70, 10).Value = 0  # CRP!J70: 34000 -> 0
$ws.Cells.Item(70, 12).Value = 0  # CRP!L70: 34000 -> 0
$ws.Cells.Item(70, 14).ClearContents()  # CRP!N70: -34630 -> (empty)
$ws.Cells.Item(73, 8).Value = 0  # CRP!H73: 34000 -> 0
$ws.Cells.Item(73, 10).Value = 0  # CRP!J73: 34000 -> 0
$ws.Cells.Item(73, 12).Value = 0  # CRP!L73: 34000 -> 0
$ws.Cells.Item(73, 14).ClearContents()  # CRP!N73: -36184 -> (empty)
$ws.Cells.Item(134, 8).Value = 2231.1143  # CRP!H134: 2076.5 -> 2231.1143
$ws.Cells.Item(134, 9).Value = 1537.0358  # CRP!I134: 1540 -> 1537.0358
$ws.Cells.Item(134, 10).Value = 5007.4287  # CRP!J134: 4759 -> 5007.4287
$ws.Cells.Item(134, 11).Value = 4611.107400000001  # CRP!K134: 4620 -> 4611.107400000001
$ws.Cells.Item(134, 12).Value = 15022.2861  # CRP!L134: 14277 -> 15022.2861
$ws.Cells.Item(134, 13).Value = -2076.107400000001  # CRP!M134: -2085 -> -2076.107400000001
$ws.Cells.Item(134, 14).Value = -20092.2861  # CRP!N134: -19347 -> -20092.2861
$ws.Cells.Item(135, 8).Value = 69000  # CRP!H135: 0 -> 69000
$ws.Cells.Item(135, 10).Value = 69000  # CRP!J135: 0 -> 69000
$ws.Cells.Item(135, 12).Value = 69000  # CRP!L135: 0 -> 69000
$ws.Cells.Item(135, 14).Value = -79140  # CRP!N135: None -> -79140

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 266.25  # CUL!H23: 324.33334 -> 266.25
$ws.Cells.Item(23, 9).Value = 50  # CUL!I23: 200 -> 50
$ws.Cells.Item(23, 10).Value = 297.14285  # CUL!J23: 349.2 -> 297.14285
$ws.Cells.Item(23, 11).Value = 150  # CUL!K23: 600 -> 150
$ws.Cells.Item(23, 12).Value = 891.4285500000001  # CUL!L23: 1047.6 -> 891.4285500000001
$ws.Cells.Item(23, 13).Value = 85  # CUL!M23: -365 -> 85
$ws.Cells.Item(23, 14).Value = -1361.42855  # CUL!N23: -1517.6 -> -1361.42855
$ws.Cells.Item(86, 8).Value = 87.25  # CUL!H86: 87 -> 87.25
$ws.Cells.Item(86, 10).Value = 92  # CUL!J86: 95 -> 92
$ws.Cells.Item(86, 12).Value = 276  # CUL!L86: 285 -> 276
$ws.Cells.Item(86, 14).Value = -2648  # CUL!N86: -2657 -> -2648
$ws.Cells.Item(89, 8).Value = 87.25  # CUL!H89: 87 -> 87.25
$ws.Cells.Item(89, 10).Value = 92  # CUL!J89: 95 -> 92
$ws.Cells.Item(89, 12).Value = 828  # CUL!L89: 855 -> 828
$ws.Cells.Item(89, 14).Value = -12684  # CUL!N89: -12711 -> -12684
$ws.Cells.Item(129, 8).Value = 1984.6  # CUL!H129: 1981.5 -> 1984.6
$ws.Cells.Item(129, 10).Value = 1998.3334  # CUL!J129: 1999 -> 1998.3334
$ws.Cells.Item(129, 12).Value = 5995.0002  # CUL!L129: 5997 -> 5995.0002
$ws.Cells.Item(129, 14).Value = -15995.0002  # CUL!N129: -15997 -> -15995.0002
$ws.Cells.Item(141, 8).Value = 8557.143  # CUL!H141: 11456.25 -> 8557.143
$ws.Cells.Item(141, 9).Value = 5980  # CUL!I141: 9780 -> 5980
$ws.Cells.Item(141, 10).Value = 15000  # CUL!J141: 14250 -> 15000
$ws.Cells.Item(141, 11).Value = 17940  # CUL!K141: 29340 -> 17940
$ws.Cells.Item(141, 12).Value = 45000  # CUL!L141: 42750 -> 45000
$ws.Cells.Item(141, 13).Value = -12760  # CUL!M141: -24160 -> -12760
$ws.Cells.Item(141, 14).Value = -55360  # CUL!N141: -53110 -> -55360

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 8500400  # GSM!H3: 11333333 -> 8500400
$ws.Cells.Item(3, 10).Value = 5000800  # GSM!J3: 10000000 -> 5000800
$ws.Cells.Item(3, 12).Value = 5000800  # GSM!L3: 10000000 -> 5000800
$ws.Cells.Item(3, 14).Value = -5001032  # GSM!N3: -10000232 -> -5001032
$ws.Cells.Item(102, 8).Value = 2895  # GSM!H102: 2930.4614 -> 2895
$ws.Cells.Item(102, 9).Value = 2851.7856  # GSM!I102: 2883 -> 2851.7856
$ws.Cells.Item(102, 11).Value = 2851.7856  # GSM!K102: 2883 -> 2851.7856
$ws.Cells.Item(102, 13).Value = -1229.7856  # GSM!M102: -1261 -> -1229.7856
$ws.Cells.Item(126, 8).Value = 6711.857  # GSM!H126: 6497.875 -> 6711.857
$ws.Cells.Item(126, 9).Value = 7497.25  # GSM!I126: 6997.8 -> 7497.25
$ws.Cells.Item(126, 11).Value = 22491.75  # GSM!K126: 20993.4 -> 22491.75
$ws.Cells.Item(126, 13).Value = -20021.75  # GSM!M126: -18523.4 -> -20021.75
$ws.Cells.Item(140, 8).Value = 63461.54  # GSM!H140: 236000 -> 63461.54
$ws.Cells.Item(140, 9).Value = 60000  # GSM!I140: 0 -> 60000
$ws.Cells.Item(140, 10).Value = 82500  # GSM!J140: 236000 -> 82500
$ws.Cells.Item(140, 11).Value = 60000  # GSM!K140: 0 -> 60000
$ws.Cells.Item(140, 12).Value = 82500  # GSM!L140: 236000 -> 82500
$ws.Cells.Item(140, 13).Value = -54820  # GSM!M140: None -> -54820
$ws.Cells.Item(140, 14).Value = -92860  # GSM!N140: -246360 -> -92860

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 4399.4  # LTW!H17: 5374.75 -> 4399.4
$ws.Cells.Item(17, 9).Value = 3999.25  # LTW!I17: 4999.5 -> 3999.25
$ws.Cells.Item(17, 10).Value = 6000  # LTW!J17: 5750 -> 6000
$ws.Cells.Item(17, 11).Value = 3999.25  # LTW!K17: 4999.5 -> 3999.25
$ws.Cells.Item(17, 12).Value = 6000  # LTW!L17: 5750 -> 6000
$ws.Cells.Item(17, 13).Value = -3829.25  # LTW!M17: -4829.5 -> -3829.25
$ws.Cells.Item(17, 14).Value = -6340  # LTW!N17: -6090 -> -6340
$ws.Cells.Item(18, 8).Value = 26889.445  # LTW!H18: 483273.62 -> 26889.445
$ws.Cells.Item(18, 9).Value = 12502.5  # LTW!I18: 33003.332 -> 12502.5
$ws.Cells.Item(18, 10).Value = 31000  # LTW!J18: 652125 -> 31000
$ws.Cells.Item(18, 11).Value = 12502.5  # LTW!K18: 33003.332 -> 12502.5
$ws.Cells.Item(18, 12).Value = 31000  # LTW!L18: 652125 -> 31000
$ws.Cells.Item(18, 13).Value = -12330.5  # LTW!M18: -32831.332 -> -12330.5
$ws.Cells.Item(18, 14).Value = -31344  # LTW!N18: -652469 -> -31344
$ws.Cells.Item(51, 8).Value = 0  # LTW!H51: 38500 -> 0
$ws.Cells.Item(51, 10).Value = 0  # LTW!J51: 38500 -> 0
$ws.Cells.Item(51, 12).Value = 0  # LTW!L51: 38500 -> 0
$ws.Cells.Item(51, 14).ClearContents()  # LTW!N51: -39456 -> (empty)
$ws.Cells.Item(82, 8).Value = 4803.615  # LTW!H82: 2159.4075 -> 4803.615
$ws.Cells.Item(82, 9).Value = 4416.3335  # LTW!I82: 1530 -> 4416.3335
$ws.Cells.Item(82, 10).Value = 4919.8  # LTW!J82: 2592.125 -> 4919.8
$ws.Cells.Item(82, 11).Value = 4416.3335  # LTW!K82: 1530 -> 4416.3335
$ws.Cells.Item(82, 12).Value = 4919.8  # LTW!L82: 2592.125 -> 4919.8
$ws.Cells.Item(82, 13).Value = -4055.3335  # LTW!M82: -1169 -> -4055.3335
$ws.Cells.Item(82, 14).Value = -5641.8  # LTW!N82: -3314.125 -> -5641.8
$ws.Cells.Item(85, 8).Value = 4803.615  # LTW!H85: 2159.4075 -> 4803.615
$ws.Cells.Item(85, 9).Value = 4416.3335  # LTW!I85: 1530 -> 4416.3335
$ws.Cells.Item(85, 10).Value = 4919.8  # LTW!J85: 2592.125 -> 4919.8
$ws.Cells.Item(85, 11).Value = 4416.3335  # LTW!K85: 1530 -> 4416.3335
$ws.Cells.Item(85, 12).Value = 4919.8  # LTW!L85: 2592.125 -> 4919.8
$ws.Cells.Item(85, 13).Value = -3168.3335  # LTW!M85: -282 -> -3168.3335
$ws.Cells.Item(85, 14).Value = -7415.8  # LTW!N85: -5088.125 -> -7415.8
$ws.Cells.Item(136, 8).Value = 3316.3635  # LTW!H136: 3072.0527 -> 3316.3635
$ws.Cells.Item(136, 9).Value = 2536.8  # LTW!I136: 2357.2666 -> 2536.8
$ws.Cells.Item(136, 11).Value = 7610.400000000001  # LTW!K136: 7071.7998 -> 7610.400000000001
$ws.Cells.Item(136, 13).Value = -5060.400000000001  # LTW!M136: -4521.7998 -> -5060.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 0  # WVR!H17: 2000 -> 0
$ws.Cells.Item(17, 9).Value = 0  # WVR!I17: 2000 -> 0
$ws.Cells.Item(17, 11).Value = 0  # WVR!K17: 2000 -> 0
$ws.Cells.Item(17, 13).ClearContents()  # WVR!M17: -1828 -> (empty)
$ws.Cells.Item(62, 8).Value = 4851.7144  # WVR!H62: 4980.2856 -> 4851.7144
$ws.Cells.Item(62, 10).Value = 4890.727  # WVR!J62: 5054.364 -> 4890.727
$ws.Cells.Item(62, 12).Value = 4890.727  # WVR!L62: 5054.364 -> 4890.727
$ws.Cells.Item(62, 14).Value = -6138.727  # WVR!N62: -6302.364 -> -6138.727
$ws.Cells.Item(65, 8).Value = 4851.7144  # WVR!H65: 4980.2856 -> 4851.7144
$ws.Cells.Item(65, 10).Value = 4890.727  # WVR!J65: 5054.364 -> 4890.727
$ws.Cells.Item(65, 12).Value = 24453.635  # WVR!L65: 25271.82 -> 24453.635
$ws.Cells.Item(65, 14).Value = -30693.635  # WVR!N65: -31511.82 -> -30693.635
$ws.Cells.Item(70, 8).Value = 31599.8  # WVR!H70: 24000 -> 31599.8
$ws.Cells.Item(70, 9).Value = 21999.5  # WVR!I70: 22000 -> 21999.5
$ws.Cells.Item(70, 10).Value = 38000  # WVR!J70: 26000 -> 38000
$ws.Cells.Item(70, 11).Value = 21999.5  # WVR!K70: 22000 -> 21999.5
$ws.Cells.Item(70, 12).Value = 38000  # WVR!L70: 26000 -> 38000
$ws.Cells.Item(70, 13).Value = -21684.5  # WVR!M70: -21685 -> -21684.5
$ws.Cells.Item(70, 14).Value = -38630  # WVR!N70: -26630 -> -38630
$ws.Cells.Item(73, 8).Value = 31599.8  # WVR!H73: 24000 -> 31599.8
$ws.Cells.Item(73, 9).Value = 21999.5  # WVR!I73: 22000 -> 21999.5
$ws.Cells.Item(73, 10).Value = 38000  # WVR!J73: 26000 -> 38000
$ws.Cells.Item(73, 11).Value = 21999.5  # WVR!K73: 22000 -> 21999.5
$ws.Cells.Item(73, 12).Value = 38000  # WVR!L73: 26000 -> 38000
$ws.Cells.Item(73, 13).Value = -20907.5  # WVR!M73: -20908 -> -20907.5
$ws.Cells.Item(73, 14).Value = -40184  # WVR!N73: -28184 -> -40184
$ws.Cells.Item(122, 8).Value = 2865  # WVR!H122: 3124.348 -> 2865
$ws.Cells.Item(122, 9).Value = 2417.4736  # WVR!I122: 2633.2104 -> 2417.4736
$ws.Cells.Item(122, 10).Value = 4565.6  # WVR!J122: 5457.25 -> 4565.6
$ws.Cells.Item(122, 11).Value = 7252.4208  # WVR!K122: 7899.6312 -> 7252.4208
$ws.Cells.Item(122, 12).Value = 13696.8  # WVR!L122: 16371.75 -> 13696.8
$ws.Cells.Item(122, 13).Value = -4802.4208  # WVR!M122: -5449.6312 -> -4802.4208
$ws.Cells.Item(122, 14).Value = -18596.8  # WVR!N122: -21271.75 -> -18596.8
$ws.Cells.Item(132, 8).Value = 3040.2083  # WVR!H132: 3509.3 -> 3040.2083
$ws.Cells.Item(132, 9).Value = 2767.1052  # WVR!I132: 3155.375 -> 2767.1052
$ws.Cells.Item(132, 10).Value = 4078  # WVR!J132: 4925 -> 4078
$ws.Cells.Item(132, 11).Value = 8301.3156  # WVR!K132: 9466.125 -> 8301.3156
$ws.Cells.Item(132, 12).Value = 12234  # WVR!L132: 14775 -> 12234
$ws.Cells.Item(132, 13).Value = -5771.3156  # WVR!M132: -6936.125 -> -5771.3156
$ws.Cells.Item(132, 14).Value = -17294  # WVR!N132: -19835 -> -17294
